$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Paragraph 1 ("... template :"): append a fully empty run (<w:r/>) right
# after the trailing run that holds NBSP + ":" but before the paragraph
# mark. This is the fix described by the commit message ("Fixed some empty
# run when bookmarks are used.") - M2Doc was leaving a stray empty run in
# place right before the paragraph that starts the bookmarked "test" call.
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs(1)
$p1End = $p1.Range.End

# The final two characters of the paragraph (NBSP + ":") are exactly the
# content of the last run; replace that exact span with itself plus a new,
# completely empty run so nothing else about the existing run is altered.
$lastRun = $d.Range($p1End - 3, $p1End - 1)

$lastRunXml = "<pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'>" +
    "<pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'>" +
    "<pkg:xmlData>" +
    "<w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>" +
    "<w:body><w:p>" +
    "<w:r w:rsidRPr='00DC5685'><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>&#160;:</w:t></w:r>" +
    "<w:r/>" +
    "</w:p></w:body></w:document>" +
    "</pkg:xmlData></pkg:part></pkg:package>"

$lastRun.InsertXML($lastRunXml)
